# "adding area to all discharge files"
#
# Adds cross-sectional Area columns next to the existing Discharge (Q)
# columns: per-segment Area (G), running Atotal (H), and a small
# "totals" summary pair (J: Atotal, K: Qtotal) mirroring H2/F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("G1").Value   = "Area"
$ws.Range("H1").Value   = "Atotal"
$ws.Range("J1").Value   = "Atotal"
$ws.Range("K1").Value   = "Qtotal"

# --- Row 2: first segment uses 0 as the starting depth ------------------
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# --- Row 3: plain (non-shared) formula, like E3 --------------------------
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# --- Rows 4-15: shared formula, analogous to the E5:E8 group below -------
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Small summary block next to the data table --------------------------
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Match the author's final on-screen selection -------------------------
$ws.Range("J2:K2").Select() | Out-Null

$wb.Save()
